# Refresh crypto price ("D") and 1h volume-change ("E") columns with the
# latest values from the scheduled scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.114.90'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.563.75'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '581.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.107'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.356'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.54'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').Value = '3.022.22'
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = '63.035.68'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000145'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.78%  '
$ws.Range('D17').Value = '2.544.74'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.38'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '342.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('E23').Value = '  +0.38%  '
$ws.Range('D24').Value = '2.679.11'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  +2.36%  '
$ws.Range('E26').Value = '  +0.82%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.44'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.87'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.20%  '
$ws.Range('E31').Value = '  +6.04%  '
$ws.Range('D32').Value = '0.0₃0823'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '176.81'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.45%  '
$ws.Range('E34').Value = '  +0.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '428.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.90%  '
$ws.Range('E36').Value = '  +1.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.22'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.94%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.45'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').Value = '  +0.51%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.66'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '151.97'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.26%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0550'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.51%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.606'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.21%  '
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0242'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.41'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.58%  '
$ws.Range('E51').Value = '  -3.81%  '
